$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.590.34'
$ws.Range('E2').Value = '  +2.42%  '
$ws.Range('D3').Value = '3.131.08'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '617.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  +0.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.403'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.39%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').Value = '3.128.33'
$ws.Range('E10').Value = '  +30.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.760'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.26%  '
$ws.Range('D14').Value = '93.272.44'
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.75'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '3.723.23'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '3.189.53'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('E20').Value = '  -0.98%  '
$ws.Range('B21').Value = 'PEPE'
$ws.Range('C21').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000208'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.49%  '
$ws.Range('B22').Value = 'Polkadot'
$ws.Range('C22').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.83'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '450.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.61'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.86'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').Value = '3.297.40'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.137'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +11.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.230'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.171'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.30'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('E34').Value = '  +6.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.159'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '26.37'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '496.35'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.81%  '
$ws.Range('E39').Value = '  -0.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.85'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.31%  '
$ws.Range('E41').Value = '  -3.21%  '
$ws.Range('E42').Value = '  -1.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.31%  '
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '162.28'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.93'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.694'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0337'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.14%  '
